$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.819.24"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "2.191.61"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "291.43"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "86.50"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").Value = "0.508"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.464"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").Value = "29.91"
$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("E11").Value = "  +6.29%  "

$ws.Range("D12").Value = "0.0776"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").Value = "2.528.77"
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "13.63"
$ws.Range("E16").Value = "  -3.05%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.195.80"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").Value = "39.710.33"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").Value = "11.17"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").Value = "5.70"
$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").Value = "235.94"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -1.31%  "

$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").Value = "22.91"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("E29").Value = "  -7.90%  "

$ws.Range("D30").Value = "9.12"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").Value = "155.81"
$ws.Range("E31").Value = "  +2.29%  "

$ws.Range("D32").Value = "31.09"
$ws.Range("E32").Value = "  -6.33%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").Value = "0.0703"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").Value = "0.0969"
$ws.Range("E39").Value = "  -3.23%  "

$ws.Range("D40").Value = "15.04"
$ws.Range("E40").Value = "  -7.26%  "

$ws.Range("E41").Value = "  -3.37%  "

$ws.Range("D42").Value = "2.121.19"
$ws.Range("E42").Value = "  +2.69%  "

$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -2.57%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0266"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.07"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").Value = "9.66"
$ws.Range("E46").Value = "  -2.24%  "

$ws.Range("D47").Value = "17.03"
$ws.Range("E47").Value = "  -5.07%  "

$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  +2.22%  "

$ws.Range("D49").Value = "2.401.30"
$ws.Range("E49").Value = "  -1.52%  "

$ws.Range("D50").Value = "1.46"
$ws.Range("E50").Value = "  +1.30%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "87.76"
$ws.Range("E51").Value = "  -1.54%  "
